$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.319.23"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.21%  "

$ws.Range("D3").Value = "'2.615.90"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.11%  "

$ws.Range("E4").Value = "  -0.04%  "

$ws.Range("D5").Value = "'551.99"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.68%  "

$ws.Range("D6").Value = "'154.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.38%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").Value = "'0.594"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.07%  "

$ws.Range("E9").Value = "  -3.08%  "

$ws.Range("D11").Value = "'5.45"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.12%  "

$ws.Range("D12").Value = "'0.365"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.38%  "

$ws.Range("D13").Value = "'3.076.95"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.23%  "

$ws.Range("D14").Value = "'25.66"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.61%  "

$ws.Range("D15").Value = "'62.223.00"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.12%  "

$ws.Range("E16").Value = "  -3.00%  "

$ws.Range("D17").Value = "'2.615.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.24%  "

$ws.Range("D18").Value = "'11.64"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.92%  "

$ws.Range("D19").Value = "'4.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -3.24%  "

$ws.Range("D20").Value = "'340.86"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.36%  "

$ws.Range("E21").Value = "  -6.43%  "

$ws.Range("D22").Value = "'1.00"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.38%  "

$ws.Range("D23").Value = "'0.497"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.09%  "

$ws.Range("D24").Value = "'62.99"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.51%  "

$ws.Range("D25").Value = "'0.167"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.34%  "

$ws.Range("E26").Value = "  +0.04%  "

$ws.Range("D27").Value = "'8.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -2.95%  "

$ws.Range("B28").Value = "PEPE"
$ws.Range("C28").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D28").Value = "'0.0₃0827"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -7.51%  "

$ws.Range("D29").Value = "'7.16"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.45%  "

$ws.Range("B30").Value = "Fetch.AI"
$ws.Range("C30").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D30").Value = "'1.34"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.73%  "

$ws.Range("E31").Value = "  -3.08%  "

$ws.Range("B32").Value = "Monero"
$ws.Range("C32").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D32").Value = "'160.26"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.86%  "

$ws.Range("B33").Value = "USDe"
$ws.Range("C33").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D33").Value = "'0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("D34").Value = "'19.25"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.36%  "

$ws.Range("E35").Value = "  -3.64%  "

$ws.Range("E36").Value = "  -4.67%  "

$ws.Range("E37").Value = "  -3.13%  "

$ws.Range("D38").Value = "'336.78"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.20%  "

$ws.Range("D39").Value = "'6.12"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.33%  "

$ws.Range("D40").Value = "'0.890"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.69%  "

$ws.Range("D41").Value = "'37.68"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.86%  "

$ws.Range("E42").Value = "  -3.97%  "

$ws.Range("D43").Value = "'0.997"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.15%  "

$ws.Range("B44").Value = "InjectiveProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D44").Value = "'20.41"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.75%  "

$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.610"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.80%  "

$ws.Range("D46").Value = "'2.129.98"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.22%  "

$ws.Range("D47").Value = "'10.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.11%  "

$ws.Range("E48").Value = "  -5.32%  "

$ws.Range("E49").Value = "  -5.41%  "

$ws.Range("D50").Value = "'0.0966"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.41%  "

$ws.Range("D51").Value = "'0.0239"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.19%  "
